$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.189873417721519
$ws.Range("C2").Value = 0.5854430379746836
$ws.Range("J2").Value = 0.006329113924050633
$ws.Range("P2").Value = 0.1487341772151899
$ws.Range("S2").Value = 0.06962025316455696
$ws.Range("B3").Value = 0.0160427807486631
$ws.Range("C3").Value = 0.0267379679144385
$ws.Range("J3").Value = 0.0160427807486631
$ws.Range("P3").Value = 0.7807486631016043
$ws.Range("S3").Value = 0.160427807486631
$ws.Range("J4").Value = 0.08695652173913043
$ws.Range("P4").Value = 0.5434782608695652
$ws.Range("S4").Value = 0.3695652173913043
$ws.Range("B6").Value = 0.08
$ws.Range("D6").Value = 0.01
$ws.Range("F6").Value = 0.06
$ws.Range("J6").Value = 0.26
$ws.Range("O6").Value = 0.005
$ws.Range("Q6").Value = 0.185
$ws.Range("R6").Value = 0.05
$ws.Range("S6").Value = 0.35
$ws.Range("B7").Value = 0.1421800947867299
$ws.Range("D7").Value = 0.03791469194312796
$ws.Range("F7").Value = 0.05213270142180094
$ws.Range("J7").Value = 0.1421800947867299
$ws.Range("O7").Value = 0.01421800947867299
$ws.Range("Q7").Value = 0.1374407582938389
$ws.Range("R7").Value = 0.06635071090047394
$ws.Range("S7").Value = 0.4075829383886256
$ws.Range("B8").Value = 0.1240875912408759
$ws.Range("D8").Value = 0.0170316301703163
$ws.Range("E8").Value = 0.0024330900243309
$ws.Range("F8").Value = 0.0583941605839416
$ws.Range("J8").Value = 0.0827250608272506
$ws.Range("O8").Value = 0.009732360097323601
$ws.Range("Q8").Value = 0.170316301703163
$ws.Range("R8").Value = 0.09975669099756691
$ws.Range("S8").Value = 0.4355231143552312
$ws.Range("B9").Value = 0.0963855421686747
$ws.Range("D9").Value = 0.01204819277108434
$ws.Range("F9").Value = 0.05220883534136546
$ws.Range("J9").Value = 0.1244979919678715
$ws.Range("O9").Value = 0.004016064257028112
$ws.Range("Q9").Value = 0.1927710843373494
$ws.Range("R9").Value = 0.1044176706827309
$ws.Range("S9").Value = 0.4136546184738956
$ws.Range("B10").Value = 0.1018808777429467
$ws.Range("D10").Value = 0.02037617554858934
$ws.Range("F10").Value = 0.06269592476489028
$ws.Range("J10").Value = 0.1175548589341693
$ws.Range("O10").Value = 0.01018808777429467
$ws.Range("Q10").Value = 0.1896551724137931
$ws.Range("R10").Value = 0.07601880877742946
$ws.Range("S10").Value = 0.4216300940438871
$ws.Range("G11").Value = 0.14
$ws.Range("J11").Value = 0.07333333333333333
$ws.Range("K11").Value = 0.1533333333333333
$ws.Range("L11").Value = 0.6266666666666667
$ws.Range("S11").Value = 0.006666666666666667
$ws.Range("G12").Value = 0.7598039215686274
$ws.Range("J12").Value = 0.1323529411764706
$ws.Range("L12").Value = 0.06862745098039216
$ws.Range("S12").Value = 0.0392156862745098
$ws.Range("G13").Value = 0.6774193548387096
$ws.Range("J13").Value = 0.2580645161290323
$ws.Range("S13").Value = 0.06451612903225806
$ws.Range("F15").Value = 0.0128755364806867
$ws.Range("H15").Value = 0.1545064377682404
$ws.Range("I15").Value = 0.1072961373390558
$ws.Range("J15").Value = 0.3991416309012876
$ws.Range("K15").Value = 0.04291845493562232
$ws.Range("M15").Value = 0.01716738197424893
$ws.Range("O15").Value = 0.04721030042918455
$ws.Range("S15").Value = 0.2188841201716738
$ws.Range("F16").Value = 0.02314814814814815
$ws.Range("H16").Value = 0.1481481481481481
$ws.Range("I16").Value = 0.1203703703703704
$ws.Range("J16").Value = 0.4212962962962963
$ws.Range("K16").Value = 0.1435185185185185
$ws.Range("M16").Value = 0.01388888888888889
$ws.Range("O16").Value = 0.03240740740740741
$ws.Range("S16").Value = 0.09722222222222222
$ws.Range("F17").Value = 0.01886792452830189
$ws.Range("H17").Value = 0.1485849056603774
$ws.Range("I17").Value = 0.1108490566037736
$ws.Range("J17").Value = 0.4127358490566038
$ws.Range("K17").Value = 0.1108490566037736
$ws.Range("M17").Value = 0.009433962264150943
$ws.Range("O17").Value = 0.08490566037735849
$ws.Range("S17").Value = 0.1037735849056604
$ws.Range("F18").Value = 0.005319148936170213
$ws.Range("H18").Value = 0.1436170212765958
$ws.Range("I18").Value = 0.148936170212766
$ws.Range("J18").Value = 0.4095744680851064
$ws.Range("K18").Value = 0.07446808510638298
$ws.Range("M18").Value = 0.01063829787234043
$ws.Range("O18").Value = 0.101063829787234
$ws.Range("S18").Value = 0.1063829787234043
$ws.Range("F19").Value = 0.009923664122137405
$ws.Range("H19").Value = 0.1923664122137405
$ws.Range("I19").Value = 0.09465648854961832
$ws.Range("J19").Value = 0.3801526717557252
$ws.Range("K19").Value = 0.1145038167938931
$ws.Range("M19").Value = 0.01450381679389313
$ws.Range("N19").Value = 0.002290076335877863
$ws.Range("O19").Value = 0.08778625954198473
$ws.Range("S19").Value = 0.1038167938931298
